$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "AAA"
$ws.Range("B6").Value = "MrPlFV/OrrrDHZeiO+tz1YB6vkINsvdtc4n1hyUNjVs="
$ws.Range("C6").Value = "GjNm+Zhdefr7W9LtefsSDw=="

$ws.Range("A7").Value = "k"
$ws.Range("B7").Value = "NWH3CJSReh6JsrIXloaA+5EPvPgJyt1yt4P5kgF4AkU="
$ws.Range("C7").Value = "SGl5a0zLko4BN+bAmwBWOA=="
